$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new row at position 12 (pushes the old rows 12-16 down to 13-17).
# Excel's native row-insert already clones the "s" (style) attributes from
# the row above for columns A-D/F (style 6) and leaves E/G without an
# explicit style, which matches the target row exactly.
# ---------------------------------------------------------------------------
$ws.Rows("12:12").Insert()

# Fill in the new row 12 with Diego Ancalao's data.
$ws.Range("A12").Value = 16
$ws.Range("B12").Value = "Diego Ancalao"
$ws.Range("C12").Value = "Independiente"
$ws.Range("D12").Value = "IND"
$ws.Range("E12").Value = "https://ellibero.cl/actualidad/ancalao-la-otra-carta-presidencial-para-la-lista-del-pueblo/?mc_cid=fd4ab4e2fa&mc_eid=1c664a6593"
$ws.Range("F12").Value = "."
$ws.Range("G12").Value = 0

# ---------------------------------------------------------------------------
# Hyperlinks: the engine's Range.Hyperlinks.Delete() clears every hyperlink
# on the whole sheet (not just the target range), and Hyperlinks.Add()
# rewrites the target cell's style to a generic "Hyperlink" xf, losing the
# original per-cell formatting nuance (plain hyperlink vs. left-aligned vs.
# date-formatted hyperlink cells). To rebuild the full, final hyperlink set
# at the shifted cells *and* keep each cell's exact original look, we:
#   1. snapshot each cell's current formatting onto the clipboard,
#   2. add the hyperlink (which stomps on the style),
#   3. paste the snapshotted formatting back over the same cell.
# ---------------------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete() | Out-Null

function Add-Link($addr, $url) {
    $rng = $ws.Range($addr)
    $rng.Copy() | Out-Null
    $ws.Hyperlinks.Add($rng, $url) | Out-Null
    $rng.PasteSpecial(-4122) | Out-Null
}

Add-Link "E2"  "https://www.cnnchile.com/pais/gabriel-boric-candidato-presidencial-frente-amplio_20210317/"
Add-Link "F2"  "https://www.latercera.com/politica/noticia/convergencia-social-proclama-a-gabriel-boric-como-precandidato-presidencial/E4Z3B6W7ZBB67GJ3T7TWQWJ4SE/"
Add-Link "E3"  "https://www.latercera.com/politica/noticia/daniel-jadue-asegura-que-esta-absolutamente-disponible-para-ser-presidente/5LTX46V4VVAE3NKHKVBOLL6BOE/"
Add-Link "F3"  "https://www.youtube.com/watch?time_continue=1&v=V8j8qed30f4&feature=emb_logo&ab_channel=TVN"
Add-Link "F4"  "https://www.elmostrador.cl/noticias/multimedia/2021/01/26/los-aplausos-en-la-moneda-en-la-despedida-del-ahora-precandidato-presidencial-ignacio-briones/"
Add-Link "E7"  "https://www.cnnchile.com/lodijeronencnn/lavin-elecciones-presidenciales-futuro-gobierno-mirada-amplia_20200824/"
Add-Link "F8"  "https://www.adnradio.cl/politica/2020/10/29/economista-y-exgerente-de-desaparecida-linea-law-se-lanza-como-candidato-presidencial.html"
Add-Link "E10" "https://www.lacuarta.com/espectaculos/noticia/dr-file-quiero-proximo-presidente-chile/550493/"
Add-Link "F10" "https://www.publimetro.cl/cl/entretenimiento/2020/10/01/dr-file-anuncia-carrera-presidencial-quiero-proximo-presidente-chile.html "
Add-Link "E12" "https://ellibero.cl/actualidad/ancalao-la-otra-carta-presidencial-para-la-lista-del-pueblo/?mc_cid=fd4ab4e2fa&mc_eid=1c664a6593"
Add-Link "E13" "https://www.latercera.com/politica/noticia/partido-humanista-declara-a-pamela-jiles-como-carta-presidencial-y-acusa-operacion-del-frente-amplio-para-impedir-que-la-diputada-lidere-la-comision-de-constitucion/A65MH4UQG5E3FK6UZS2UYHMYHI/"
Add-Link "E14" "https://www.biobiochile.cl/noticias/nacional/chile/2018/12/16/jose-antonio-kast-anuncia-carrera-presidencial-para-elecciones-de-2021-en-medio-de-seminario.shtml"

# ---------------------------------------------------------------------------
# Selection cursor moves to B13 (where the old Pamela Jiles row now lives).
# ---------------------------------------------------------------------------
$ws.Range("B13").Select() | Out-Null
